$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 153.875
$ws.Range("I39").Value = 153.875
$ws.Range("K39").Value = 461.625
$ws.Range("M39").Value = -165.625
$ws.Range("H43").Value = 4248
$ws.Range("I43").Value = 2871.75
$ws.Range("J43").Value = 5165.5
$ws.Range("K43").Value = 2871.75
$ws.Range("L43").Value = 5165.5
$ws.Range("M43").Value = -2802.75
$ws.Range("N43").Value = -5303.5
$ws.Range("H58").Value = 205.33333
$ws.Range("I58").Value = 199.5
$ws.Range("J58").Value = 217
$ws.Range("K58").Value = 598.5
$ws.Range("L58").Value = 651
$ws.Range("M58").Value = -448.5
$ws.Range("N58").Value = -951
$ws.Range("H104").Value = 300
$ws.Range("I104").Value = 300
$ws.Range("K104").Value = 900
$ws.Range("M104").Value = 847
$ws.Range("H113").Value = 8142.2104
$ws.Range("I113").Value = 10223.556
$ws.Range("J113").Value = 6269
$ws.Range("K113").Value = 10223.556
$ws.Range("L113").Value = 6269
$ws.Range("M113").Value = -6969.556
$ws.Range("N113").Value = -12777
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H132").Value = 1390.8889
$ws.Range("I132").Value = 1390.8889
$ws.Range("K132").Value = 4172.6667
$ws.Range("M132").Value = -1642.6667
$ws.Range("H137").Value = 3734
$ws.Range("I137").Value = 3734
$ws.Range("K137").Value = 11202
$ws.Range("M137").Value = -8652
$ws.Range("H138").Value = 1261.75
$ws.Range("J138").Value = 2199.3333
$ws.Range("L138").Value = 6597.999899999999
$ws.Range("N138").Value = -16877.9999
$ws.Range("H141").Value = 3600
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2550.8462
$ws.Range("I45").Value = 1073.4
$ws.Range("J45").Value = 3474.25
$ws.Range("K45").Value = 1073.4
$ws.Range("L45").Value = 3474.25
$ws.Range("M45").Value = -696.4000000000001
$ws.Range("N45").Value = -4228.25
$ws.Range("H61").Value = 1394.4546
$ws.Range("I61").Value = 1384.1
$ws.Range("K61").Value = 1384.1
$ws.Range("M61").Value = -1172.1
$ws.Range("H110").Value = 366.66666
$ws.Range("I110").Value = 175
$ws.Range("J110").Value = 750
$ws.Range("K110").Value = 175
$ws.Range("L110").Value = 750
$ws.Range("M110").Value = 1870
$ws.Range("N110").Value = -4840
$ws.Range("H136").Value = 1394.4546
$ws.Range("I136").Value = 1384.1
$ws.Range("K136").Value = 4152.299999999999
$ws.Range("M136").Value = -1602.299999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 3167199.8
$ws.Range("I7").Value = 3800539.8
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 3800539.8
$ws.Range("L7").Value = 500
$ws.Range("M7").Value = -3800426.8
$ws.Range("N7").Value = -726
$ws.Range("H10").Value = 5302
$ws.Range("I10").Value = 600
$ws.Range("J10").Value = 7653
$ws.Range("K10").Value = 600
$ws.Range("L10").Value = 7653
$ws.Range("M10").Value = -460
$ws.Range("N10").Value = -7933
$ws.Range("H86").Value = 3909.8
$ws.Range("I86").Value = 750
$ws.Range("K86").Value = 750
$ws.Range("M86").Value = 373
$ws.Range("H89").Value = 3909.8
$ws.Range("I89").Value = 750
$ws.Range("K89").Value = 3750
$ws.Range("M89").Value = 1866
$ws.Range("H99").Value = 1037.8948
$ws.Range("I99").Value = 1045.0555
$ws.Range("K99").Value = 1045.0555
$ws.Range("M99").Value = 452.9445000000001
$ws.Range("H134").Value = 1500
$ws.Range("I134").Value = 1500
$ws.Range("K134").Value = 4500
$ws.Range("M134").Value = -1965

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2500
$ws.Range("I4").Value = 1666.6666
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 1666.6666
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = -1554.6666
$ws.Range("N4").Value = -5224
$ws.Range("H58").Value = 3119.5
$ws.Range("I58").Value = 2212
$ws.Range("K58").Value = 2212
$ws.Range("M58").Value = -2009
$ws.Range("H60").Value = 16600
$ws.Range("I60").Value = 3000
$ws.Range("K60").Value = 3000
$ws.Range("M60").Value = -2489
$ws.Range("H105").Value = 2796.3125
$ws.Range("I105").Value = 934.3333
$ws.Range("J105").Value = 5190.2856
$ws.Range("K105").Value = 934.3333
$ws.Range("L105").Value = 5190.2856
$ws.Range("M105").Value = 812.6667
$ws.Range("N105").Value = -8684.285599999999
$ws.Range("H136").Value = 3119.5
$ws.Range("I136").Value = 2212
$ws.Range("K136").Value = 6636
$ws.Range("M136").Value = -4086

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1280.4445
$ws.Range("I7").Value = 2050
$ws.Range("J7").Value = 318.5
$ws.Range("K7").Value = 6150
$ws.Range("L7").Value = 955.5
$ws.Range("M7").Value = -6038
$ws.Range("N7").Value = -1179.5
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H34").Value = 425
$ws.Range("I34").Value = 342.85715
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 1028.57145
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -944.5714499999999
$ws.Range("N34").Value = -3168
$ws.Range("H139").Value = 3552.5
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 449
$ws.Range("I13").Value = 315
$ws.Range("K13").Value = 315
$ws.Range("M13").Value = -176
$ws.Range("H17").Value = 1604.9
$ws.Range("I17").Value = 40
$ws.Range("J17").Value = 1778.7778
$ws.Range("K17").Value = 40
$ws.Range("L17").Value = 1778.7778
$ws.Range("M17").Value = 128
$ws.Range("N17").Value = -2114.7778
$ws.Range("H80").Value = 6668.5
$ws.Range("I80").Value = 2005
$ws.Range("J80").Value = 7601.2
$ws.Range("K80").Value = 2005
$ws.Range("L80").Value = 7601.2
$ws.Range("M80").Value = -1007
$ws.Range("N80").Value = -9597.200000000001
$ws.Range("H83").Value = 6668.5
$ws.Range("I83").Value = 2005
$ws.Range("J83").Value = 7601.2
$ws.Range("K83").Value = 10025
$ws.Range("L83").Value = 38006
$ws.Range("M83").Value = -5033
$ws.Range("N83").Value = -47990
$ws.Range("H113").Value = 5034.143
$ws.Range("I113").Value = 3848
$ws.Range("J113").Value = 7999.5
$ws.Range("K113").Value = 3848
$ws.Range("L113").Value = 7999.5
$ws.Range("M113").Value = -1678
$ws.Range("N113").Value = -12339.5
$ws.Range("H132").Value = 4377.364
$ws.Range("I132").Value = 4593
$ws.Range("J132").Value = 3802.3333
$ws.Range("K132").Value = 13779
$ws.Range("L132").Value = 11406.9999
$ws.Range("M132").Value = -11249
$ws.Range("N132").Value = -16466.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1542.25
$ws.Range("I7").Value = 1441.6666
$ws.Range("J7").Value = 1844
$ws.Range("K7").Value = 1441.6666
$ws.Range("L7").Value = 1844
$ws.Range("M7").Value = -1329.6666
$ws.Range("N7").Value = -2068
$ws.Range("H22").Value = 2319.3
$ws.Range("I22").Value = 1100
$ws.Range("J22").Value = 3132.1667
$ws.Range("K22").Value = 1100
$ws.Range("L22").Value = 3132.1667
$ws.Range("M22").Value = -805
$ws.Range("N22").Value = -3722.1667
$ws.Range("H27").Value = 2319.3
$ws.Range("I27").Value = 1100
$ws.Range("J27").Value = 3132.1667
$ws.Range("K27").Value = 1100
$ws.Range("L27").Value = 3132.1667
$ws.Range("M27").Value = -993
$ws.Range("N27").Value = -3346.1667
$ws.Range("H40").Value = 4212
$ws.Range("I40").Value = 4247.3335
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 4247.3335
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -4111.3335
$ws.Range("N40").Value = -4272
$ws.Range("H47").Value = 18500
$ws.Range("J47").Value = 18500
$ws.Range("L47").Value = 18500
$ws.Range("N47").Value = -19480
$ws.Range("H52").Value = 18500
$ws.Range("J52").Value = 18500
$ws.Range("L52").Value = 18500
$ws.Range("N52").Value = -18966
$ws.Range("H126").Value = 1542.25
$ws.Range("I126").Value = 1441.6666
$ws.Range("J126").Value = 1844
$ws.Range("K126").Value = 4324.9998
$ws.Range("L126").Value = 5532
$ws.Range("M126").Value = -1854.9998
$ws.Range("N126").Value = -10472

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 6333.3335
$ws.Range("J5").Value = 6333.3335
$ws.Range("L5").Value = 6333.3335
$ws.Range("N5").Value = -6557.3335
$ws.Range("H20").Value = 10000000
$ws.Range("J20").Value = 10000000
$ws.Range("L20").Value = 10000000
$ws.Range("N20").Value = -10000480
$ws.Range("I22").Value = 1800
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1800
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1507
$ws.Range("N22").ClearContents()
$ws.Range("H107").Value = 1351
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1351
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 4053
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -7893
$ws.Range("H117").Value = 35000
$ws.Range("J117").Value = 35000
$ws.Range("L117").Value = 35000
$ws.Range("N117").Value = -44178
$ws.Range("H124").Value = 35000
$ws.Range("J124").Value = 35000
$ws.Range("L124").Value = 35000
$ws.Range("N124").Value = -44820
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -70120
$ws.Range("H136").Value = 2700
$ws.Range("I136").Value = 2700
$ws.Range("K136").Value = 8100
$ws.Range("M136").Value = -5550
